# Weekly fruit/vegetable price update:
# Insert a new record (row) above the existing row 158, shifting the
# remaining historical rows down by one, and populate the new row with
# the latest "Coliflor" price observation for Macroferia Regional de
# Talca (Maule).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 158..168 down to 159..169, inserting a fresh blank row at 158.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A158").Value = 5
$ws.Range("B158").Value = "Macroferia Regional de Talca"
$ws.Range("C158").Value = "Maule"
$ws.Range("D158").Value = 44516
$ws.Range("E158").Value = 7
$ws.Range("F158").Value = 100112008
$ws.Range("G158").Value = "Coliflor"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 5000
$ws.Range("K158").Value = 500
$ws.Range("L158").Value = 500
$ws.Range("M158").Value = 500
$ws.Range("N158").Value = '$/unidad'
$ws.Range("O158").Value = "Región del Maule"
$ws.Range("P158").Value = 500
$ws.Range("Q158").Value = 1
$ws.Range("R158").Value = "Hortaliza"
